$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, copying the style from H1 (bold/centered/bordered header style)
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @{
  2  = @(9, 9)
  3  = @(5, 5)
  4  = @(6, 7)
  5  = @(8, 8)
  6  = @(8, 8)
  7  = @(8, 8)
  8  = @(4, 6)
  9  = @(5, 6)
  10 = @(6, 7)
  11 = @(4, 5)
  12 = @(5, 5)
  13 = @(9, 9)
  14 = @(5, 5)
  15 = @(7, 8)
  16 = @(10, 10)
  17 = @(3, 5)
  18 = @(9, 9)
  19 = @(11, 11)
  20 = @(9, 9)
  21 = @(8, 9)
  22 = @(6, 6)
  23 = @(5, 6)
  24 = @(8, 9)
  25 = @(7, 8)
  26 = @(5, 8)
  27 = @(6, 8)
  28 = @(8, 8)
  29 = @(10, 10)
  30 = @(7, 7)
  31 = @(9, 9)
  32 = @(5, 5)
  33 = @(4, 4)
  34 = @(7, 7)
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $ws.Cells.Item($row, 9).Value = $vals[0]
  $ws.Cells.Item($row, 10).Value = $vals[1]
}
